$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A16")
$cell.Formula = '=IF(ISBLANK(B16), "Düsseldorf", B16)'
$cell.Select()
